{"js": "// The document contains a paragraph holding a Word field whose field code is\n// `{ m:'doc.html'.fromHTMLURI() }` (an M2Doc template field). The edit\n// \"flattens\" that field: the field delimiters (begin/end fldChar) and the\n// instrText runs that make up the field code are replaced by plain visible\n// text runs, using literal \"{\" and \"}\" braces in place of the field\n// begin/end, while keeping the bookmark (_GoBack) and the textual content of\n// the field code runs intact.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Load each paragraph's fields so we can find the one holding our field.\nfor (const paragraph of paragraphs.items) {\n  paragraph.fields.load(\"items\");\n}\nawait context.sync();\n\n// Locate the paragraph that contains the (single) field.\nlet fieldParagraph = null;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.fields.items.length > 0) {\n    fieldParagraph = paragraph;\n    break;\n  }\n}\n\nif (fieldParagraph) {\n  // Replace the whole paragraph's content (the field) with the equivalent\n  // plain-text runs via an OOXML fragment, mirroring the target diff:\n  //   fldChar(begin) + \" \"            -> <w:t>{</w:t>\n  //   instrText \"m\"                   -> <w:t>m</w:t>\n  //   instrText \":\"                   -> <w:t>:</w:t>\n  //   instrText \"'\"                   -> <w:t>'</w:t>\n  //   instrText \"doc.html\"            -> <w:t>doc.html</w:t>\n  //   bookmarkStart/bookmarkEnd \"_GoBack\" -> unchanged\n  //   instrText \"'.fromHTMLURI()\"     -> <w:t>'.fromHTMLURI()</w:t>\n  //   \" \" + fldChar(end)              -> <w:t xml:space=\"preserve\">}</w:t>\n  const range = fieldParagraph.getRange();\n  const ooxml =\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r><w:t>{</w:t></w:r>' +\n    '<w:r><w:t>m</w:t></w:r>' +\n    '<w:r><w:t>:</w:t></w:r>' +\n    \"<w:r><w:t>'</w:t></w:r>\" +\n    '<w:r><w:t>doc.html</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    \"<w:r><w:t>'.fromHTMLURI()</w:t></w:r>\" +\n    '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n  range.insertOoxml(ooxml, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# The document contains a paragraph holding a Word field whose field code is\n# `{ m:'doc.html'.fromHTMLURI() }` (an M2Doc template field). The edit\n# \"flattens\" that field: the field delimiters (begin/end fldChar) and the\n# instrText runs that make up the field code are replaced by plain visible\n# text runs, using literal \"{\" and \"}\" braces in place of the field\n# begin/end, while keeping the bookmark (_GoBack) and the textual content of\n# the field code runs intact.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the (single) field, instead of\n# hard-coding a paragraph index.\n$targetParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($candidate.Range.Fields.Count -gt 0) {\n        $targetParagraph = $candidate\n        break\n    }\n}\n\nif ($targetParagraph -ne $null) {\n    # Replace the whole paragraph's content (the field) with the equivalent\n    # plain-text runs via a WordOpenXML fragment, mirroring the target diff:\n    #   fldChar(begin) + \" \"                -> <w:t>{</w:t>\n    #   instrText \"m\"                       -> <w:t>m</w:t>\n    #   instrText \":\"                       -> <w:t>:</w:t>\n    #   instrText \"'\"                       -> <w:t>'</w:t>\n    #   instrText \"doc.html\"                -> <w:t>doc.html</w:t>\n    #   bookmarkStart/bookmarkEnd \"_GoBack\" -> unchanged\n    #   instrText \"'.fromHTMLURI()\"         -> <w:t>'.fromHTMLURI()</w:t>\n    #   \" \" + fldChar(end)                  -> <w:t xml:space=\"preserve\">}</w:t>\n    $apostrophe = \"'\"\n    $xml = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n        '<w:p>' +\n        '<w:r><w:t>{</w:t></w:r>' +\n        '<w:r><w:t>m</w:t></w:r>' +\n        '<w:r><w:t>:</w:t></w:r>' +\n        '<w:r><w:t>' + $apostrophe + '</w:t></w:r>' +\n        '<w:r><w:t>doc.html</w:t></w:r>' +\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n        '<w:bookmarkEnd w:id=\"0\"/>' +\n        '<w:r><w:t>' + $apostrophe + '.fromHTMLURI()</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n        '</w:p>' +\n        '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n\n    $targetParagraph.Range.InsertXML($xml)\n}\n"}
